$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - column F holds "想去人数" (want-to-go count)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 84
$wsExhibit.Range("F4").Value = 7206
$wsExhibit.Range("F5").Value = 263
$wsExhibit.Range("F6").Value = 415
$wsExhibit.Range("F7").Value = 3687
$wsExhibit.Range("F8").Value = 306
$wsExhibit.Range("F9").Value = 530
$wsExhibit.Range("F10").Value = 268
$wsExhibit.Range("F11").Value = 601
$wsExhibit.Range("F12").Value = 90

# Sheet "全部类型" (All types) - same underlying rows, offset by the extra
# "演出" (performance) rows interleaved, still column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 84
$wsAll.Range("F6").Value = 7206
$wsAll.Range("F8").Value = 263
$wsAll.Range("F9").Value = 415
$wsAll.Range("F10").Value = 3687
$wsAll.Range("F11").Value = 306
$wsAll.Range("F12").Value = 530
$wsAll.Range("F13").Value = 268
$wsAll.Range("F14").Value = 601
$wsAll.Range("F15").Value = 90
